# New crime data collected - weekly CompStat update for the 68th Precinct.
# Updates the report week/volume header text and refreshes the crime-count
# table (weekly, 28-day, YTD, 2-year columns) for rows 16-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: bump the report volume/number and the covered week dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  23"
$ws.Range("C9").Value = "Report Covering the Week  6/5/2023  Through  6/11/2023"

# ---------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = 20
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = 56.521739130434
$ws.Range("L16").Value = 44
$ws.Range("M16").Value = -21.739130434782
$ws.Range("N16").Value = -86.567164179104

# ---------------------------------------------------------------------
# Row 17 (Fel. Assault) - D17/E17 flip from numbers to the text markers
# used elsewhere in the sheet for "no data" ('0'/'***.*').
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = "'0"
$ws.Range("A17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "'***.*"
$ws.Range("A17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 350
$ws.Range("I17").Value = 66
$ws.Range("K17").Value = 57.142857142857
$ws.Range("L17").Value = 65
$ws.Range("M17").Value = 43.478260869565
$ws.Range("N17").Value = -43.103448275862

# ---------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 40
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 33
$ws.Range("K18").Value = 57.575757575757
$ws.Range("L18").Value = 18.181818181818
$ws.Range("M18").Value = -50.943396226415
$ws.Range("N18").Value = -89.233954451345

# ---------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 17
$ws.Range("E19").Value = -35.294117647058
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = -30
$ws.Range("I19").Value = 194
$ws.Range("J19").Value = 246
$ws.Range("K19").Value = -21.138211382113
$ws.Range("L19").Value = 31.081081081081
$ws.Range("M19").Value = 46.969696969697
$ws.Range("N19").Value = -6.280193236714

# ---------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 51
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = 6.25
$ws.Range("L20").Value = 88.888888888888
$ws.Range("M20").Value = -22.727272727272
$ws.Range("N20").Value = -94.308035714285

# ---------------------------------------------------------------------
# Row 21 (TOTAL, bold row)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -13.043478260869
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = -9.333333333333
$ws.Range("I21").Value = 407
$ws.Range("J21").Value = 400
$ws.Range("K21").Value = 1.75
$ws.Range("L21").Value = 40.830449826989
$ws.Range("M21").Value = 0.992555831265
$ws.Range("N21").Value = -79.485887096774

# ---------------------------------------------------------------------
# Row 22 (Transit) - C22/F22 flip from the text '0' marker to real
# numbers, while D22/E22 flip the other way, to the text markers.
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 1
$ws.Range("G22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = "'***.*"
$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("G22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = -55.555555555555
$ws.Range("L22").Value = 33.333333333333
$ws.Range("M22").Value = 0

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 7.407407407407
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = 4.950495049504
$ws.Range("I24").Value = 685
$ws.Range("J24").Value = 783
$ws.Range("K24").Value = -12.515964240102
$ws.Range("L24").Value = 83.155080213903
$ws.Range("M24").Value = 31.226053639846

# ---------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 94.117647058823
$ws.Range("I25").Value = 170
$ws.Range("J25").Value = 140
$ws.Range("K25").Value = 21.428571428571
$ws.Range("L25").Value = 63.461538461538
$ws.Range("M25").Value = 8.280254777070

# ---------------------------------------------------------------------
# Row 27 (Other Sex Crimes)
# ---------------------------------------------------------------------
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = -52.380952380952
